# "Add files via upload" - append the 0.7.2 version-history row to Sheet1
# and move the active selection to C24, matching the new upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row of version history data (row 26) ------------------------------
$ws.Range("A26").Value = "0.7.2"
$ws.Range("B26").Value = "AUTOMATA CELULAR - copia (36)"
$ws.Range("C26").Value = "-UI: Delete rows according to working functionality.`n-Document every function.`n_OPTIONAL: Make it possible to reduce Niches on mutations.`n-When there is no data, the program crashes."
$ws.Range("D26").Value = "-Reworked Aggrupation.`n-Reworked Flexibility.`n-Lambda can now be negative."
$ws.Range("E26").Value = "Python 3.6.1"
$ws.Range("F26").Value = "Qt version: 5.6.2`nSIP version: 4.18`nPyQt version: 5.6"

# Row grows to fit the wrapped multi-line text, same as the other
# "changelog" rows above it.
$ws.Rows.Item(26).RowHeight = 57.6

# --- Selection moves to C24 -------------------------------------------------
$ws.Range("C24").Select() | Out-Null
